# Frederick County residential-density zone descriptions were reworded:
# the trailing sentence ("The following residential density districts and
# maximum densities are hereby established.") was dropped from the shared
# paragraph used in column D (Zone General Description) for the six
# residential rows (R-1, R-3, R-5, R-8, R-12, R-16).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Frederick County")

$newText = "The purpose of the residential density districts is to promote healthful and convenient distribution of population with sufficient densities to maintain a high standard of physical design and community service. Residential density districts will conform to the County Comprehensive Plan and will be located within areas identified for residential development. The districts, as a group, are intended to provide for a variety of dwelling types and densities and to offer housing choices at various economic levels. It is further the intent to establish various densities of residential developments in order to efficiently and effectively provide for necessary public services and facilities."

# Row 3 (R-1) gets the clean, trailing-space-free text.
$ws.Range("D3").Value = $newText

# Rows 4-8 (R-3, R-5, R-8, R-12, R-16) get the same text but with a
# trailing space, matching how the author's edit propagated.
$ws.Range("D4").Value = $newText + " "
$ws.Range("D5").Value = $newText + " "
$ws.Range("D6").Value = $newText + " "
$ws.Range("D7").Value = $newText + " "
$ws.Range("D8").Value = $newText + " "

# The author's last interaction on this sheet left the selection on T5
# (sheet stays inactive - Washington County remains the active tab).
$ws.Range("T5").Select()

$wsWash = $wb.Worksheets.Item("Washington County")
$wsWash.Range("D30").Select()
